# Downtime Tracker data-folder refresh
# ----------------------------------------------------------------------------
# The upstream report changed its Downtime column from "days" to "hours"
# (values multiplied by 24, with previously-unset/zero rows floored to a
# minimum half-day of 12h) and corrected a handful of mis-dated rows that had
# rolled into the wrong reporting months. Two previously-blank trailing rows
# also received their missing Status/Downtime/Resolution data.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Downtime Tracker")

# --- Column A (Date) corrections: rows had drifted into Oct-Dec 2024,
#     should read Jan-Feb 2024 ---
$ws.Cells.Item(28, 1).Value = 45309
$ws.Cells.Item(29, 1).Value = 45312
$ws.Cells.Item(30, 1).Value = 45323
$ws.Cells.Item(31, 1).Value = 45330
$ws.Cells.Item(32, 1).Value = 45336

# --- Column H (Downtime / Resolution Time), converted from days to hours ---
$h = @{
    2  = 48;  3  = 12;  4  = 12;  5  = 264; 6  = 456; 7  = 48;  8  = 24;
    9  = 48;  10 = 12;  11 = 24;  12 = 24;  13 = 48;  14 = 24;  15 = 12;
    16 = 12;  17 = 48;  18 = 12;  19 = 48;  20 = 24;  21 = 12;  22 = 48;
    23 = 48;  24 = 24;  25 = 24;  26 = 12;  27 = 12;  28 = 12;  29 = 12;
    30 = 12;  31 = 12;  32 = 48;  33 = 12;  34 = 12;  35 = 12;  36 = 12;
    37 = 12;  38 = 12;  39 = 12;  40 = 12;  41 = 12;  42 = 12;  43 = 48;
    44 = 24;  45 = 48;  46 = 24;  47 = 12;  48 = 24;  49 = 12
}
foreach ($row in $h.Keys | Sort-Object) {
    $ws.Cells.Item($row, 8).Value = $h[$row]
}

# --- Row 49 was missing its Status / Permanent Resolution text; fill it in
#     and match the formatting already used by the rest of the table (the
#     row had been left with the blank "unused row" style) ---
$ws.Range("G47").Copy() | Out-Null
$ws.Range("G49").PasteSpecial(-4122) | Out-Null
$ws.Range("I47").Copy() | Out-Null
$ws.Range("I49").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(49, 7).Value = "Resolved"
$ws.Cells.Item(49, 9).Value = "Currently Unclear"

# --- View state: scrolled down to row 24, with A28 selected ---
$ws.Activate()
$ws.Range("A28").Select() | Out-Null
